$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update column F (想去人数) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 271
$ws1.Range("F3").Value = 74
$ws1.Range("F5").Value = 7107
$ws1.Range("F6").Value = 5499
$ws1.Range("F7").Value = 453
$ws1.Range("F9").Value = 9
$ws1.Range("F12").Value = 149

# Sheet "全部类型" (fourth sheet) - same updates, plus row 14 instead of row 12
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 271
$ws4.Range("F3").Value = 74
$ws4.Range("F5").Value = 7107
$ws4.Range("F6").Value = 5499
$ws4.Range("F7").Value = 453
$ws4.Range("F9").Value = 9
$ws4.Range("F14").Value = 149
